# Updated cryptos list on Tue Oct  3 14:38:01 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) cells for rows
# 2-51 of the crypto table with newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price strings look like plain numbers (e.g. "214.10",
# "1.20") even though the column stores them as literal text (to keep
# trailing zeros / exact source formatting, matching the other Price cells
# like "27.457.64" that aren't valid numbers at all). Mark those cells as
# Text before writing so Excel doesn't silently coerce them into numeric
# values and drop the trailing zero. ClearFormats() afterwards drops the
# temporary number-format override again so the cell style stays the
# default, unstyled one - only the underlying value stays text.
$textRefs = @("D5", "D6", "D9", "D11", "D15", "D16", "D18", "D20", "D23", "D25", "D26", "D27", "D31", "D37", "D42", "D43", "D48", "D51")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.457.64'
$ws.Range("E2").Value = '  -3.21%  '

$ws.Range("D3").Value = '1.653.54'
$ws.Range("E3").Value = '  -3.63%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '214.10'
$ws.Range("E5").Value = '  -2.13%  '

$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  -2.28%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("E8").Value = '  +1.20%  '

$ws.Range("D9").Value = '0.262'
$ws.Range("E9").Value = '  -1.65%  '

$ws.Range("E10").Value = '  -2.42%  '

$ws.Range("D11").Value = '0.0875'
$ws.Range("E11").Value = '  -1.78%  '

$ws.Range("D12").Value = '1.887.23'
$ws.Range("E12").Value = '  -3.53%  '

$ws.Range("D13").Value = '1.657.30'
$ws.Range("E13").Value = '  -3.29%  '

$ws.Range("E14").Value = '  -2.24%  '

$ws.Range("D15").Value = '0.565'
$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("D16").Value = '65.69'
$ws.Range("E16").Value = '  -2.53%  '

$ws.Range("D17").Value = '27.441.96'
$ws.Range("E17").Value = '  -3.10%  '

$ws.Range("D18").Value = '236.22'
$ws.Range("E18").Value = '  -4.79%  '

$ws.Range("E19").Value = '  -2.64%  '

$ws.Range("D20").Value = '7.56'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("E22").Value = '  -3.31%  '

$ws.Range("D23").Value = '9.31'
$ws.Range("E23").Value = '  -3.18%  '

$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").Value = '145.83'
$ws.Range("E25").Value = '  -1.30%  '

$ws.Range("D26").Value = '7.18'
$ws.Range("E26").Value = '  -3.09%  '

$ws.Range("D27").Value = '16.10'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("E29").Value = '  -2.34%  '

$ws.Range("E30").Value = '  -2.65%  '

$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  -0.85%  '

$ws.Range("E32").Value = '  -2.88%  '

$ws.Range("D33").Value = '1.448.02'
$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("E34").Value = '  -4.56%  '

$ws.Range("E35").Value = '  -4.36%  '

$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").Value = '0.915'
$ws.Range("E37").Value = '  -6.02%  '

$ws.Range("E38").Value = '  -4.31%  '

$ws.Range("E39").Value = '  -2.88%  '

$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").Value = '66.26'
$ws.Range("E42").Value = '  -4.63%  '

$ws.Range("D43").Value = '5.45'
$ws.Range("E43").Value = '  -3.42%  '

$ws.Range("E44").Value = '  -2.44%  '

$ws.Range("E45").Value = '  -2.16%  '

$ws.Range("D46").Value = '1.795.98'
$ws.Range("E46").Value = '  -3.45%  '

$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").Value = '88.53'
$ws.Range("E48").Value = '  -1.69%  '

$ws.Range("E49").Value = '  -1.68%  '

$ws.Range("E50").Value = '  -2.04%  '

$ws.Range("D51").Value = '7.81'
$ws.Range("E51").Value = '  -3.56%  '

# Drop the temporary text-format override so these cells keep the default
# (unstyled) cell format, matching the rest of the column.
foreach ($r in $textRefs) {
    $ws.Range($r).ClearFormats()
}
